# Add ApDocument model test
# - Adds "created_by"/"updated_by" columns (C, D) to the "counterparties" sheet,
#   matching the shape already used on the other sheets (users/roles/roles_users/
#   company_units/employees).
# - Converts the counterparties' eik_egn (B) column from numeric to a quoted text
#   value (e.g. 111222333 -> "111222333"), matching how IDs are stored elsewhere.
# - Leaves a trailing selection on C5:D8 (the newly added block) on every sheet,
#   as left behind by the editing session.

$wb = $excel.ActiveWorkbook

$users         = $wb.Worksheets.Item("users")
$roles         = $wb.Worksheets.Item("roles")
$rolesUsers    = $wb.Worksheets.Item("roles_users")
$companyUnits  = $wb.Worksheets.Item("company_units")
$employees     = $wb.Worksheets.Item("employees")
$counterparties = $wb.Worksheets.Item("counterparties")

# --- counterparties: add created_by / updated_by columns -------------------
$counterparties.Activate()

$counterparties.Range("C1").Value = "created_by"
$counterparties.Range("D1").Value = "updated_by"

$ids = @("111222333", "111222334", "111222335", "111222336", "111222337", "111222338", "111222339")
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = 2 + $i
    $counterparties.Range("B$row").Value = '"' + $ids[$i] + '"'
    $counterparties.Range("C$row").Value = 1
    $counterparties.Range("D$row").Value = 1
}

# --- leftover selections from the editing session ---------------------------
$users.Activate()
$users.Range("C5:D8").Select()

$roles.Activate()
$roles.Range("B6").Select()

$rolesUsers.Activate()
$rolesUsers.Range("N17").Select()

$companyUnits.Activate()
$companyUnits.Range("B4").Select()

$employees.Activate()
$employees.Range("E1").Select()

$counterparties.Activate()
$counterparties.Range("C5:D8").Select()
